# "control dia 27 febrero" - add two new time-tracking entries to the
# "Fabio" sheet (sheet2) and make that sheet the active tab, mirroring
# the selection that Excel leaves behind after the user's last edit.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- New row 7: 2015-02-27, "crear portafolio" ---
$ws2.Range("A7").Value = 42062
$ws2.Range("B7").Value = "crear portafolio"
$ws2.Range("C7").Value = 4.5
$ws2.Range("F7").Value = 0

# --- New row 8: 2015-02-24, "terminar logo" ---
$ws2.Range("A8").Value = 42059
$ws2.Range("B8").Value = "terminar logo"
$ws2.Range("C8").Value = 4.5
$ws2.Range("F8").Value = 100

# Give the new date cells the same date formatting as the existing ones
# (copy/paste-special formats only, so no new style entries are created).
$ws2.Range("A6").Copy() | Out-Null
$ws2.Range("A7:A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Select F9 on the "Fabio" sheet and make it the active tab/sheet, just
# like it was left in the saved workbook.
$ws2.Range("F9").Select() | Out-Null

Write-Host "done"
